$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A13:A15 with the new customer ids
$ws.Range("A13").Value = "cus_KQ6aaMpTMkZP9V"
$ws.Range("A14").Value = "cus_KQ6ayxL6jsGJ8c"
$ws.Range("A15").Value = "cus_KQ6aP84xmpHq5P"

# Remove rows 16 through 21 which held the old, now-tested customer ids
$ws.Range("A16:A21").EntireRow.Delete()

# Update selection to reflect the newly added/verified range
$ws.Range("A13:A15").Select()
